$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the affected D/E cells to Text format so percentage/decimal strings
# are preserved exactly as text (matching original inlineStr storage),
# rather than being auto-converted to numbers by Excel.
$priceVolRange = $ws.Range("D2:E47")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "277.03"
$ws.Range("E2").Value = "0.67%"
$ws.Range("D3").Value = "27.26"
$ws.Range("E3").Value = "0.48%"
$ws.Range("D4").Value = "4.841"
$ws.Range("E4").Value = "1.78%"
$ws.Range("E5").Value = "0.48%"
$ws.Range("D6").Value = "7.025"
$ws.Range("E6").Value = "1.30%"
$ws.Range("D7").Value = "1.362"
$ws.Range("E7").Value = "0.50%"
$ws.Range("D8").Value = "0.8884"
$ws.Range("E8").Value = "1.34%"
$ws.Range("D9").Value = "0.1511"
$ws.Range("E9").Value = "0.13%"
$ws.Range("D10").Value = "0.05482"
$ws.Range("E10").Value = "8.26%"
$ws.Range("D11").Value = "0.07464"
$ws.Range("E11").Value = "-0.51%"
$ws.Range("D12").Value = "0.02905"
$ws.Range("E12").Value = "-0.30%"
$ws.Range("D13").Value = "0.08954"
$ws.Range("E13").Value = "-0.60%"
$ws.Range("D14").Value = "0.001576"
$ws.Range("E14").Value = "0.34%"
$ws.Range("D15").Value = "0.0006330"
$ws.Range("E15").Value = "-0.58%"
$ws.Range("D16").Value = "0.006030"
$ws.Range("E16").Value = "1.69%"
$ws.Range("D17").Value = "3.474"
$ws.Range("E17").Value = "0.76%"
$ws.Range("E18").Value = "-0.17%"
$ws.Range("E19").Value = "-1.73%"
$ws.Range("D21").Value = "0.1340"
$ws.Range("E21").Value = "1.19%"
$ws.Range("D22").Value = "3.920"
$ws.Range("E22").Value = "-0.02%"
$ws.Range("D23").Value = "0.1506"
$ws.Range("E23").Value = "9.13%"
$ws.Range("D24").Value = "0.04375"
$ws.Range("E24").Value = "-0.46%"
$ws.Range("E25").Value = "0.32%"
$ws.Range("D26").Value = "0.004247"
$ws.Range("E26").Value = "10.74%"
$ws.Range("D28").Value = "0.0001179"
$ws.Range("E28").Value = "-1.77%"
$ws.Range("E29").Value = "-14.92%"
$ws.Range("D40").Value = "0.03989"
$ws.Range("E40").Value = "-2.71%"
$ws.Range("D41").Value = "0.006693"
$ws.Range("E41").Value = "-1.41%"
$ws.Range("D42").Value = "0.1395"
$ws.Range("E42").Value = "19.16%"
$ws.Range("D43").Value = "0.002159"
$ws.Range("E43").Value = "1.78%"
$ws.Range("D44").Value = "0.01180"
$ws.Range("E44").Value = "2.32%"
$ws.Range("D45").Value = "0.00005554"
$ws.Range("E45").Value = "7.07%"
$ws.Range("E47").Value = "-19.63%"

# Restore default (Normal) style so no stray number-format style is left
# attached to these cells, matching the original workbook formatting.
$priceVolRange.Style = "Normal"
